$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.402.68'
$ws.Range('E2').Value = '  -1.97%  '
$ws.Range('D3').Value = '3.378.92'
$ws.Range('E3').Value = '  -1.95%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '567.43'
$ws.Range('E5').Value = '  -2.37%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '139.96'
$ws.Range('E6').Value = '  -6.17%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '3.379.04'
$ws.Range('E8').Value = '  -1.98%  '
$ws.Range('E9').Value = '  -0.73%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.49'
$ws.Range('E10').Value = '  -3.65%  '
$ws.Range('E11').Value = '  -3.08%  '
$ws.Range('E12').Value = '  -1.33%  '
$ws.Range('D13').Value = '3.956.10'
$ws.Range('E13').Value = '  -1.95%  '
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '27.94'
$ws.Range('E15').Value = '  +0.17%  '
$ws.Range('D16').Value = '3.380.66'
$ws.Range('E16').Value = '  -1.95%  '
$ws.Range('E17').Value = '  -3.22%  '
$ws.Range('D18').Value = '60.556.14'
$ws.Range('E18').Value = '  -1.88%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.17'
$ws.Range('E19').Value = '  -1.74%  '
$ws.Range('E20').Value = '  -4.10%  '
$ws.Range('E21').Value = '  -5.41%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '385.92'
$ws.Range('E22').Value = '  -0.76%  '
$ws.Range('E23').Value = '  -2.43%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '72.99'
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('E26').Value = '  -8.14%  '
$ws.Range('D27').Value = '3.519.49'
$ws.Range('E27').Value = '  -1.94%  '
$ws.Range('E28').Value = '  -2.09%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.36'
$ws.Range('E30').Value = '  -4.84%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.90'
$ws.Range('E31').Value = '  -4.15%  '
$ws.Range('E32').Value = '  -2.25%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.40'
$ws.Range('E33').Value = '  -8.64%  '
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '23.39'
$ws.Range('E35').Value = '  -2.81%  '
$ws.Range('D36').Value = '3.409.81'
$ws.Range('E36').Value = '  -1.79%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '168.06'
$ws.Range('E37').Value = '  +1.26%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.85'
$ws.Range('E38').Value = '  -2.44%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.96'
$ws.Range('E39').Value = '  -4.94%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.49'
$ws.Range('E40').Value = '  -5.02%  '
$ws.Range('E41').Value = '  -2.34%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '27.02'
$ws.Range('E42').Value = '  -0.40%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.778'
$ws.Range('E44').Value = '  -1.63%  '
$ws.Range('E45').Value = '  -2.48%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '41.41'
$ws.Range('E46').Value = '  -2.17%  '
$ws.Range('E47').Value = '  -1.46%  '
$ws.Range('D48').Value = '2.517.44'
$ws.Range('E48').Value = '  -3.35%  '
$ws.Range('E49').Value = '  -4.37%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '23.01'
$ws.Range('E50').Value = '  -1.03%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '6.72'
$ws.Range('E51').Value = '  -3.41%  '
